# Adds the new match row (row 23) to the india/isl/2023-2024 odds sheet,
# mirroring the existing data rows: Indice=22, Bengaluru FC vs Goa.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 23
$templateRow = 22

# Copy formatting (styles only) from the row above so the new row keeps the
# same look (bold/bordered/centered index cell, date-time formatted cell).
$ws.Cells.Item($templateRow, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122)

$ws.Cells.Item($templateRow, 5).Copy()
$ws.Cells.Item($newRow, 5).PasteSpecial(-4122)

# Column values for the new match.
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "india"
$ws.Range("C23").Value = "isl"
$ws.Range("D23").Value = "2023-2024"
$ws.Range("E23").Value = 45224.6875
$ws.Range("F23").Value = "Bengaluru FC"
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = "Goa"
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 2.53
$ws.Range("K23").Value = "21/10/2023 13:12"
$ws.Range("L23").Value = 2.54
$ws.Range("M23").Value = "25/10/2023 16:28"
$ws.Range("N23").Value = 3.35
$ws.Range("O23").Value = "21/10/2023 13:12"
$ws.Range("P23").Value = 3.68
$ws.Range("Q23").Value = "25/10/2023 16:26"
$ws.Range("R23").Value = 2.78
$ws.Range("S23").Value = "21/10/2023 13:12"
$ws.Range("T23").Value = 2.64
$ws.Range("U23").Value = "25/10/2023 16:28"
$ws.Range("V23").Value = "https://www.betexplorer.com/football/india/isl/bengaluru-fc-fc-goa/dSyKPCW3/"
